$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K3: "3X" (component count multiplier label) ---
$ws.Range("K3").Value = "3X"
$ws.Range("B3").Copy()
$ws.Range("K3").PasteSpecial(-4122)

# --- Data cells K4:K38: component counts for 3 sets (=3*G{row}) ---
$ws.Range("K4").Formula = "=3*G4"
$ws.Range("K5:K38").Formula = "=3*G5"

$ws.Range("G4").Copy()
$ws.Range("K4:K38").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update the active selection to match the saved view ---
[void]$ws.Range("B24").Select()

Write-Host "done"
